# Auto-generated: apply scheduled-runner market data refresh to all leve profit sheets.
# For each (sheet, cell) pair we set the new numeric value; cells that the refresh
# removed entirely (now blank/no profit) are cleared with $null so the cell disappears
# from the saved XML, matching the source diff exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39657
$ws.Range("J3").Value = 39657
$ws.Range("L3").Value = 39657
$ws.Range("N3").Value = -39885
$ws.Range("H40").Value = 4152
$ws.Range("I40").Value = 2725
$ws.Range("K40").Value = 2725
$ws.Range("M40").Value = -2550
$ws.Range("H58").Value = 3334
$ws.Range("J58").Value = 3334
$ws.Range("L58").Value = 10002
$ws.Range("N58").Value = -10302
$ws.Range("H95").Value = 14915.667
$ws.Range("J95").Value = 14915.667
$ws.Range("L95").Value = 14915.667
$ws.Range("N95").Value = -20407.667
$ws.Range("H98").Value = 2740.4546
$ws.Range("I98").Value = 1235
$ws.Range("J98").Value = 3995
$ws.Range("K98").Value = 1235
$ws.Range("L98").Value = 3995
$ws.Range("M98").Value = 263
$ws.Range("N98").Value = -6991
$ws.Range("H102").Value = 39657
$ws.Range("J102").Value = 39657
$ws.Range("L102").Value = 39657
$ws.Range("N102").Value = -46147
$ws.Range("H118").Value = 2688.5
$ws.Range("I118").Value = 2600
$ws.Range("K118").Value = 7800
$ws.Range("M118").Value = -6143
$ws.Range("H122").Value = 2740.4546
$ws.Range("I122").Value = 1235
$ws.Range("J122").Value = 3995
$ws.Range("K122").Value = 3705
$ws.Range("L122").Value = 11985
$ws.Range("M122").Value = -1255
$ws.Range("N122").Value = -16885
$ws.Range("H127").Value = 5804.25
$ws.Range("I127").Value = 7000
$ws.Range("K127").Value = 21000
$ws.Range("M127").Value = -16040
$ws.Range("H129").Value = 3547.6
$ws.Range("I129").Value = 1374.5
$ws.Range("J129").Value = 4996.3335
$ws.Range("K129").Value = 4123.5
$ws.Range("L129").Value = 14989.0005
$ws.Range("M129").Value = 876.5
$ws.Range("N129").Value = -24989.0005
$ws.Range("H141").Value = 5742.6665
$ws.Range("I141").Value = 5742.6665
$ws.Range("K141").Value = 17227.9995
$ws.Range("M141").Value = -12047.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 58333
$ws.Range("J24").Value = 58333
$ws.Range("L24").Value = 58333
$ws.Range("N24").Value = -59081
$ws.Range("H32").Value = 9410.5
$ws.Range("I32").Value = 9410.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9410.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9123.5
$ws.Range("N32").Value = $null
$ws.Range("H35").Value = 1774.5
$ws.Range("I35").Value = 1774.5
$ws.Range("K35").Value = 1774.5
$ws.Range("M35").Value = -1368.5
$ws.Range("H100").Value = 58333
$ws.Range("J100").Value = 58333
$ws.Range("L100").Value = 58333
$ws.Range("N100").Value = -60497
$ws.Range("H114").Value = 39800
$ws.Range("J114").Value = 39800
$ws.Range("L114").Value = 39800
$ws.Range("N114").Value = -48478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 41714
$ws.Range("J100").Value = 41714
$ws.Range("L100").Value = 41714
$ws.Range("N100").Value = -43878
$ws.Range("H103").Value = 19899.834
$ws.Range("J103").Value = 19899.834
$ws.Range("L103").Value = 19899.834
$ws.Range("N103").Value = -22243.834
$ws.Range("H107").Value = 1893
$ws.Range("I107").Value = 1190.6666
$ws.Range("K107").Value = 1190.6666
$ws.Range("M107").Value = 729.3334
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = $null
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 5000
$ws.Range("M60").Value = -4489
$ws.Range("H80").Value = 32128
$ws.Range("J80").Value = 32128
$ws.Range("L80").Value = 32128
$ws.Range("N80").Value = -34374
$ws.Range("H83").Value = 32128
$ws.Range("J83").Value = 32128
$ws.Range("L83").Value = 96384
$ws.Range("N83").Value = -107616
$ws.Range("H96").Value = 33415.668
$ws.Range("J96").Value = 33415.668
$ws.Range("L96").Value = 33415.668
$ws.Range("N96").Value = -38907.668
$ws.Range("H99").Value = 6099.9
$ws.Range("I99").Value = 5708.3335
$ws.Range("K99").Value = 5708.3335
$ws.Range("M99").Value = -4210.3335
$ws.Range("H126").Value = 6099.9
$ws.Range("I126").Value = 5708.3335
$ws.Range("K126").Value = 17125.0005
$ws.Range("M126").Value = -14655.0005
$ws.Range("H132").Value = 1710.9
$ws.Range("I132").Value = 1388.75
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 4166.25
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -1636.25
$ws.Range("N132").Value = -14058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 30.416666
$ws.Range("I38").Value = 27.8
$ws.Range("K38").Value = 83.40000000000001
$ws.Range("M38").Value = 263.6
$ws.Range("H40").Value = 116.46154
$ws.Range("I40").Value = 30.571428
$ws.Range("J40").Value = 216.66667
$ws.Range("K40").Value = 122.285712
$ws.Range("L40").Value = 866.66668
$ws.Range("M40").Value = -53.285712
$ws.Range("N40").Value = -1004.66668
$ws.Range("H92").Value = 1200
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H98").Value = 295
$ws.Range("J98").Value = 295
$ws.Range("L98").Value = 885
$ws.Range("N98").Value = -3881
$ws.Range("H131").Value = 2655.111
$ws.Range("I131").Value = 2099.4
$ws.Range("J131").Value = 3349.75
$ws.Range("K131").Value = 6298.200000000001
$ws.Range("L131").Value = 10049.25
$ws.Range("M131").Value = -1258.200000000001
$ws.Range("N131").Value = -20129.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 30328.5
$ws.Range("J101").Value = 30328.5
$ws.Range("L101").Value = 30328.5
$ws.Range("N101").Value = -36818.5
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4467.3335
$ws.Range("I82").Value = 4402
$ws.Range("J82").Value = 4500
$ws.Range("K82").Value = 4402
$ws.Range("L82").Value = 4500
$ws.Range("M82").Value = -4041
$ws.Range("N82").Value = -5222
$ws.Range("H85").Value = 4467.3335
$ws.Range("I85").Value = 4402
$ws.Range("J85").Value = 4500
$ws.Range("K85").Value = 4402
$ws.Range("L85").Value = 4500
$ws.Range("M85").Value = -3154
$ws.Range("N85").Value = -6996
$ws.Range("H95").Value = 18000
$ws.Range("J95").Value = 18000
$ws.Range("L95").Value = 18000
$ws.Range("N95").Value = -23492
$ws.Range("H110").Value = 52500
$ws.Range("J110").Value = 52500
$ws.Range("L110").Value = 52500
$ws.Range("N110").Value = -60680
$ws.Range("H136").Value = 3336
$ws.Range("I136").Value = 3336
$ws.Range("K136").Value = 10008
$ws.Range("M136").Value = -7458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
$ws.Range("H92").Value = 46249.25
$ws.Range("J92").Value = 46249.25
$ws.Range("L92").Value = 46249.25
$ws.Range("N92").Value = -51241.25
$ws.Range("H97").Value = 58499.5
$ws.Range("J97").Value = 58499.5
$ws.Range("L97").Value = 58499.5
$ws.Range("N97").Value = -60481.5
$ws.Range("H98").Value = 33333
$ws.Range("J98").Value = 33333
$ws.Range("L98").Value = 33333
$ws.Range("N98").Value = -39323
$ws.Range("H103").Value = 30602
$ws.Range("J103").Value = 30602
$ws.Range("L103").Value = 30602
$ws.Range("N103").Value = -32946
$ws.Range("H113").Value = 1766.8235
$ws.Range("I113").Value = 912
$ws.Range("K113").Value = 2736
$ws.Range("M113").Value = -566

Write-Output "Applied 219 cell updates across 8 sheets"